# Auto-generated Excel COM-interop script applying market-data refresh
# to the "Chocobo_Profits" Leve tables across all job sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR). Updates currentAveragePrice* columns (H/I/J) and
# recomputed LevePrice*/LeveProfit* columns (K/L/M/N) per row.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 394.2857
$ws.Range("I5").Value = 265
$ws.Range("J5").Value = 566.6667
$ws.Range("K5").Value = 265
$ws.Range("L5").Value = 566.6667
$ws.Range("M5").Value = -150
$ws.Range("N5").Value = -796.6667

# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 836.3200000000001
$ws.Range("J17").Value = 702.8293
$ws.Range("L17").Value = 2108.4879
$ws.Range("N17").Value = -2444.4879

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 1897.6923
$ws.Range("I43").Value = 1056.4
$ws.Range("J43").Value = 2423.5
$ws.Range("K43").Value = 1056.4
$ws.Range("L43").Value = 2423.5
$ws.Range("M43").Value = -987.4000000000001
$ws.Range("N43").Value = -2561.5

# Row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 1319.7255
$ws.Range("I112").Value = 774.1429000000001
$ws.Range("J112").Value = 1406.5227
$ws.Range("K112").Value = 2322.4287
$ws.Range("L112").Value = 4219.5681
$ws.Range("M112").Value = -1214.4287
$ws.Range("N112").Value = -6435.5681

# Row 115 (Leve Item ID 27957)
$ws.Range("H115").Value = 984.25
$ws.Range("I115").Value = 984.25
$ws.Range("K115").Value = 2952.75
$ws.Range("M115").Value = -1385.75

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 392067.78
$ws.Range("I116").Value = 771850.1
$ws.Range("J116").Value = 12285.385
$ws.Range("K116").Value = 771850.1
$ws.Range("L116").Value = 12285.385
$ws.Range("M116").Value = -768408.1
$ws.Range("N116").Value = -19169.385

# Row 118 (Leve Item ID 27958)
$ws.Range("H118").Value = 785.8823
$ws.Range("J118").Value = 945.7778
$ws.Range("L118").Value = 2837.3334
$ws.Range("N118").Value = -6151.3334

# Row 127 (Leve Item ID 36114)
$ws.Range("H127").Value = 1626.1082
$ws.Range("J127").Value = 1866.9642
$ws.Range("L127").Value = 5600.892599999999
$ws.Range("N127").Value = -15520.8926

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2314.79
$ws.Range("I138").Value = 790.4643
$ws.Range("J138").Value = 2907.5833
$ws.Range("K138").Value = 2371.3929
$ws.Range("L138").Value = 8722.749899999999
$ws.Range("M138").Value = 2768.6071
$ws.Range("N138").Value = -19002.7499

# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 5548.9775
$ws.Range("I141").Value = 5586.1396
$ws.Range("J141").Value = 4750
$ws.Range("K141").Value = 16758.4188
$ws.Range("L141").Value = 14250
$ws.Range("M141").Value = -11578.4188
$ws.Range("N141").Value = -24610

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 997.4286
$ws.Range("I2").Value = 888.06665
$ws.Range("J2").Value = 1270.8334
$ws.Range("K2").Value = 888.06665
$ws.Range("L2").Value = 1270.8334
$ws.Range("M2").Value = -775.06665
$ws.Range("N2").Value = -1496.8334

# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 5281.5967
$ws.Range("I32").Value = 4391.365
$ws.Range("J32").Value = 9910.799999999999
$ws.Range("K32").Value = 4391.365
$ws.Range("L32").Value = 9910.799999999999
$ws.Range("M32").Value = -4104.365
$ws.Range("N32").Value = -10484.8

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 1576.25
$ws.Range("I45").Value = 1638.6666
$ws.Range("J45").Value = 1513.8334
$ws.Range("K45").Value = 1638.6666
$ws.Range("L45").Value = 1513.8334
$ws.Range("M45").Value = -1261.6666
$ws.Range("N45").Value = -2267.8334

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 4373.12
$ws.Range("I74").Value = 4066.4348
$ws.Range("J74").Value = 7900
$ws.Range("K74").Value = 4066.4348
$ws.Range("L74").Value = 7900
$ws.Range("M74").Value = -3192.4348
$ws.Range("N74").Value = -9648

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 4373.12
$ws.Range("I77").Value = 4066.4348
$ws.Range("J77").Value = 7900
$ws.Range("K77").Value = 20332.174
$ws.Range("L77").Value = 39500
$ws.Range("M77").Value = -15964.174
$ws.Range("N77").Value = -48236

# Row 109 (Leve Item ID 25646)
$ws.Range("H109").Value = 32500
$ws.Range("J109").Value = 32500
$ws.Range("L109").Value = 32500
$ws.Range("N109").Value = -35274

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 997.4286
$ws.Range("I116").Value = 888.06665
$ws.Range("J116").Value = 1270.8334
$ws.Range("K116").Value = 888.06665
$ws.Range("L116").Value = 1270.8334
$ws.Range("M116").Value = 1405.93335
$ws.Range("N116").Value = -5858.8334

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 1668.125
$ws.Range("J122").Value = 15000
$ws.Range("L122").Value = 45000
$ws.Range("N122").Value = -49900

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 1632.9767
$ws.Range("I132").Value = 803.58826
$ws.Range("J132").Value = 4766.222
$ws.Range("K132").Value = 2410.76478
$ws.Range("L132").Value = 14298.666
$ws.Range("M132").Value = 119.23522
$ws.Range("N132").Value = -19358.666

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 997.4286
$ws.Range("I3").Value = 888.06665
$ws.Range("J3").Value = 1270.8334
$ws.Range("K3").Value = 888.06665
$ws.Range("L3").Value = 1270.8334
$ws.Range("M3").Value = -774.06665
$ws.Range("N3").Value = -1498.8334

# Row 28 (Leve Item ID 19546)
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = ""

# Row 42 (Leve Item ID 22903)
$ws.Range("H42").Value = 69800
$ws.Range("J42").Value = 69800
$ws.Range("L42").Value = 69800
$ws.Range("N42").Value = -70456

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 1835.1578
$ws.Range("I134").Value = 1222.1852
$ws.Range("J134").Value = 3339.7273
$ws.Range("K134").Value = 3666.5556
$ws.Range("L134").Value = 10019.1819
$ws.Range("M134").Value = -1131.5556
$ws.Range("N134").Value = -15089.1819

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 12502275
$ws.Range("I31").Value = 988.46155
$ws.Range("J31").Value = 35718950
$ws.Range("K31").Value = 988.46155
$ws.Range("L31").Value = 35718950
$ws.Range("M31").Value = -693.46155
$ws.Range("N31").Value = -35719540

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 12502275
$ws.Range("I34").Value = 988.46155
$ws.Range("J34").Value = 35718950
$ws.Range("K34").Value = 988.46155
$ws.Range("L34").Value = 35718950
$ws.Range("M34").Value = -786.46155
$ws.Range("N34").Value = -35719354

# Row 87 (Leve Item ID 11929)
$ws.Range("H87").Value = 21085.715
$ws.Range("J87").Value = 21085.715
$ws.Range("L87").Value = 21085.715
$ws.Range("N87").Value = -23457.715

# Row 90 (Leve Item ID 11929)
$ws.Range("H90").Value = 21085.715
$ws.Range("J90").Value = 21085.715
$ws.Range("L90").Value = 63257.145
$ws.Range("N90").Value = -75113.145

# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 544.0714
$ws.Range("J107").Value = 1083
$ws.Range("L107").Value = 1083
$ws.Range("N107").Value = -4923

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2629.7693
$ws.Range("I132").Value = 2007.9773
$ws.Range("J132").Value = 6049.625
$ws.Range("K132").Value = 6023.9319
$ws.Range("L132").Value = 18148.875
$ws.Range("M132").Value = -3493.9319
$ws.Range("N132").Value = -23208.875

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 9804832
$ws.Range("J131").Value = 964.617
$ws.Range("L131").Value = 2893.851
$ws.Range("N131").Value = -12973.851

# Row 133 (Leve Item ID 44073)
$ws.Range("H133").Value = 2471.9048
$ws.Range("I133").Value = 2428.75
$ws.Range("J133").Value = 2498.4614
$ws.Range("K133").Value = 7286.25
$ws.Range("L133").Value = 7495.3842
$ws.Range("M133").Value = -2226.25
$ws.Range("N133").Value = -17615.3842

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 58.22222
$ws.Range("I2").Value = 40.666668
$ws.Range("J2").Value = 93.333336
$ws.Range("K2").Value = 40.666668
$ws.Range("L2").Value = 93.333336
$ws.Range("M2").Value = 72.333332
$ws.Range("N2").Value = -319.333336

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 3719.9092
$ws.Range("I122").Value = 2376.625
$ws.Range("J122").Value = 7302
$ws.Range("K122").Value = 7129.875
$ws.Range("L122").Value = 21906
$ws.Range("M122").Value = -4679.875
$ws.Range("N122").Value = -26806

# Row 123 (Leve Item ID 34150)
$ws.Range("H123").Value = 11052.929
$ws.Range("J123").Value = 11052.929
$ws.Range("L123").Value = 11052.929
$ws.Range("N123").Value = -15952.929

# Row 140 (Leve Item ID 42458)
$ws.Range("H140").Value = 38727.242
$ws.Range("J140").Value = 38727.242
$ws.Range("L140").Value = 38727.242
$ws.Range("N140").Value = -49087.242

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 8197.046
$ws.Range("I40").Value = 8399.666999999999
$ws.Range("J40").Value = 8056.769
$ws.Range("K40").Value = 8399.666999999999
$ws.Range("L40").Value = 8056.769
$ws.Range("M40").Value = -8263.666999999999
$ws.Range("N40").Value = -8328.769

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2203.4666
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 1773.2307
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 1773.2307
$ws.Range("M46").Value = -4812
$ws.Range("N46").Value = -2149.2307

# Row 69 (Leve Item ID 10671)
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""

# Row 72 (Leve Item ID 10671)
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 3970.5
$ws.Range("I122").Value = 2660.077
$ws.Range("J122").Value = 5280.923
$ws.Range("K122").Value = 7980.231000000001
$ws.Range("L122").Value = 15842.769
$ws.Range("M122").Value = -5530.231000000001
$ws.Range("N122").Value = -20742.769

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 1976.8975
$ws.Range("I136").Value = 1045.2812
$ws.Range("J136").Value = 6235.7144
$ws.Range("K136").Value = 3135.8436
$ws.Range("L136").Value = 18707.1432
$ws.Range("M136").Value = -585.8435999999997
$ws.Range("N136").Value = -23807.1432

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 6804187
$ws.Range("I132").Value = 922.39026
$ws.Range("J132").Value = 41670916
$ws.Range("K132").Value = 2767.17078
$ws.Range("L132").Value = 125012748
$ws.Range("M132").Value = -237.1707799999999
$ws.Range("N132").Value = -125017808

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 2662.8572
$ws.Range("I136").Value = 804.1667
$ws.Range("J136").Value = 6718.1816
$ws.Range("K136").Value = 2412.5001
$ws.Range("L136").Value = 20154.5448
$ws.Range("M136").Value = 137.4998999999998
$ws.Range("N136").Value = -25254.5448
